# "got loops kinda working"
# Switch the plan-table template placeholders from the old `{table:...}`
# loop syntax to the new `{#...}` loop syntax, and leave the selection
# where the author last left it (G11) before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan")

# Row 7 holds the per-row template placeholders for the "planData" loop.
$ws.Range("B7").Value = "{#planData.name}"
$ws.Range("C7").Value = "{#planData.role.name}"

# Move the active selection to match the author's saved cursor position.
$ws.Range("G11").Select()
